# Generate Report for Handback
# Updates status of the e75c870a-... file from "Ready for handoff" to
# "Handed back: in sync with en-US" across the Overview, zh-cn and de-de
# sheets, and records the new handback datetimes.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the e75c870a-... file; B = zh-cn status, C = de-de status
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# zh-cn sheet: row 3 is the e75c870a-... file; C = Status, H = Latest Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-03-18 16:54:51"

# de-de sheet: row 3 is the e75c870a-... file; C = Status, H = Latest Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-03-18 16:55:08"
